# Apply the change described by the diff:
# 1. Update the "Date" metadata value on the Metadata sheet.
# 2. Add a new row to the Elements sheet describing CarteProfessionnel.ExerciceProfessionnel.
# 3. Widen a few columns to fit the new, longer content.

$wb = $excel.ActiveWorkbook

# --- 1. Update Date value on Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# --- 2. Add new row 9 to Elements sheet ---
$els = $wb.Worksheets.Item("Elements")

$els.Range("A9").Value = "CarteProfessionnel.ExerciceProfessionnel"
$els.Range("B9").Value = "CarteProfessionnel.ExerciceProfessionnel"
$els.Range("D9").Value = ""
$els.Range("F9").Value = "1"
$els.Range("G9").Value = "1"
$els.Range("H9").Value = ""
$els.Range("I9").Value = ""
$els.Range("J9").Value = ""
$els.Range("K9").Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/ExerciceProfessionnel`n"
$els.Range("L9").Value = "Lien vers la classe ExerciceProfessionnel"
$els.Range("M9").Value = "Lien vers la classe ExerciceProfessionnel"
$els.Range("P9").Value = ""
$els.Range("R9").Value = ""
$els.Range("S9").Value = ""
$els.Range("T9").Value = ""
$els.Range("U9").Value = ""
$els.Range("V9").Value = ""
$els.Range("W9").Value = ""
$els.Range("X9").Value = ""
$els.Range("Y9").Value = ""
$els.Range("Z9").Value = ""
$els.Range("AA9").Value = ""
$els.Range("AB9").Value = ""
$els.Range("AC9").Value = ""
$els.Range("AD9").Value = ""
$els.Range("AE9").Value = ""
$els.Range("AF9").Value = "CarteProfessionnel.ExerciceProfessionnel"
$els.Range("AG9").Value = "1"
$els.Range("AH9").Value = "1"
$els.Range("AI9").Value = ""
$els.Range("AJ9").Value = ""
